$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 413.33334
$ws.Range("I4").Value = 397.14285
$ws.Range("K4").Value = 397.14285
$ws.Range("M4").Value = -283.14285
$ws.Range("H19").Value = 5984.4736
$ws.Range("J19").Value = 6214.067
$ws.Range("L19").Value = 6214.067
$ws.Range("N19").Value = -6564.067
$ws.Range("H32").Value = 16503.15
$ws.Range("J32").Value = 19958.166
$ws.Range("L32").Value = 19958.166
$ws.Range("N32").Value = -20610.166
$ws.Range("H33").Value = 498.78262
$ws.Range("I33").Value = 447.53845
$ws.Range("J33").Value = 565.4
$ws.Range("K33").Value = 447.53845
$ws.Range("L33").Value = 565.4
$ws.Range("M33").Value = -218.53845
$ws.Range("N33").Value = -1023.4
$ws.Range("H43").Value = 13857
$ws.Range("J43").Value = 9999.666999999999
$ws.Range("L43").Value = 9999.666999999999
$ws.Range("N43").Value = -10137.667
$ws.Range("H64").Value = 23054.273
$ws.Range("I64").Value = 25955.223
$ws.Range("K64").Value = 25955.223
$ws.Range("M64").Value = -25707.223
$ws.Range("H67").Value = 23054.273
$ws.Range("I67").Value = 25955.223
$ws.Range("K67").Value = 25955.223
$ws.Range("M67").Value = -25097.223
$ws.Range("H116").Value = 16923.5
$ws.Range("J116").Value = 9862.4
$ws.Range("L116").Value = 9862.4
$ws.Range("N116").Value = -16746.4
$ws.Range("H132").Value = 2816129.5
$ws.Range("I132").Value = 3263935.8
$ws.Range("K132").Value = 9791807.399999999
$ws.Range("M132").Value = -9789277.399999999
$ws.Range("H137").Value = 16694.182
$ws.Range("I137").Value = 29723
$ws.Range("K137").Value = 89169
$ws.Range("M137").Value = -86619
$ws.Range("H141").Value = 2712.25
$ws.Range("I141").Value = 2712.25
$ws.Range("K141").Value = 8136.75
$ws.Range("M141").Value = -2956.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19755.215
$ws.Range("J32").Value = 2115
$ws.Range("L32").Value = 2115
$ws.Range("N32").Value = -2689
$ws.Range("H45").Value = 3618.2632
$ws.Range("I45").Value = 2081.3845
$ws.Range("K45").Value = 2081.3845
$ws.Range("M45").Value = -1704.3845
$ws.Range("H63").Value = 2455.25
$ws.Range("I63").Value = 1933.375
$ws.Range("K63").Value = 1933.375
$ws.Range("M63").Value = -1247.375
$ws.Range("H66").Value = 2455.25
$ws.Range("I66").Value = 1933.375
$ws.Range("K66").Value = 9666.875
$ws.Range("M66").Value = -6234.875
$ws.Range("H88").Value = 6950
$ws.Range("J88").Value = 7939.3335
$ws.Range("L88").Value = 7939.3335
$ws.Range("N88").Value = -8751.333500000001
$ws.Range("H91").Value = 6950
$ws.Range("J91").Value = 7939.3335
$ws.Range("L91").Value = 7939.3335
$ws.Range("N91").Value = -10747.3335
$ws.Range("H110").Value = 886.43335
$ws.Range("I110").Value = 663.8421
$ws.Range("J110").Value = 1270.909
$ws.Range("K110").Value = 663.8421
$ws.Range("L110").Value = 1270.909
$ws.Range("M110").Value = 1381.1579
$ws.Range("N110").Value = -5360.909
$ws.Range("H122").Value = 1413.4073
$ws.Range("I122").Value = 1327.8125
$ws.Range("K122").Value = 3983.4375
$ws.Range("M122").Value = -1533.4375
$ws.Range("H132").Value = 1591.1
$ws.Range("I132").Value = 1090.1904
$ws.Range("J132").Value = 2144.7368
$ws.Range("K132").Value = 3270.5712
$ws.Range("L132").Value = 6434.2104
$ws.Range("M132").Value = -740.5711999999999
$ws.Range("N132").Value = -11494.2104
$ws.Range("H137").Value = 111490.555
$ws.Range("J137").Value = 119801.875
$ws.Range("L137").Value = 119801.875
$ws.Range("N137").Value = -130001.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 52696.688
$ws.Range("I20").Value = 75872.45
$ws.Range("J20").Value = 1710
$ws.Range("K20").Value = 75872.45
$ws.Range("L20").Value = 1710
$ws.Range("M20").Value = -75625.45
$ws.Range("N20").Value = -2204
$ws.Range("H88").Value = 74999
$ws.Range("J88").Value = 74999
$ws.Range("L88").Value = 74999
$ws.Range("N88").Value = -75811
$ws.Range("H91").Value = 74999
$ws.Range("J91").Value = 74999
$ws.Range("L91").Value = 74999
$ws.Range("N91").Value = -77807
$ws.Range("H105").Value = 1949.1
$ws.Range("I105").Value = 1642.875
$ws.Range("K105").Value = 1642.875
$ws.Range("M105").Value = 104.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 86619.86
$ws.Range("J100").Value = 86619.86
$ws.Range("L100").Value = 86619.86
$ws.Range("N100").Value = -88783.86
$ws.Range("H103").Value = 22632.5
$ws.Range("I103").Value = 22900
$ws.Range("K103").Value = 22900
$ws.Range("M103").Value = -21728
$ws.Range("H107").Value = 674.0714
$ws.Range("I107").Value = 492.8095
$ws.Range("K107").Value = 492.8095
$ws.Range("M107").Value = 1427.1905
$ws.Range("H122").Value = 1956.6666
$ws.Range("I122").Value = 1997.7576
$ws.Range("J122").Value = 1806
$ws.Range("K122").Value = 5993.2728
$ws.Range("L122").Value = 5418
$ws.Range("M122").Value = -3543.2728
$ws.Range("N122").Value = -10318
$ws.Range("H132").Value = 41036.92
$ws.Range("I132").Value = 41036.92
$ws.Range("K132").Value = 123110.76
$ws.Range("M132").Value = -120580.76

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1349.0454
$ws.Range("J92").Value = 857
$ws.Range("L92").Value = 2571
$ws.Range("N92").Value = -5067
$ws.Range("H97").Value = 1286.85
$ws.Range("J97").Value = 1448.7333
$ws.Range("L97").Value = 4346.199900000001
$ws.Range("N97").Value = -5338.199900000001
$ws.Range("H103").Value = 1083.4286
$ws.Range("I103").Value = 633.3333
$ws.Range("J103").Value = 1421
$ws.Range("K103").Value = 1899.9999
$ws.Range("L103").Value = 4263
$ws.Range("M103").Value = -1020.9999
$ws.Range("N103").Value = -6021
$ws.Range("H107").Value = 1367.9395
$ws.Range("I107").Value = 2084.5
$ws.Range("J107").Value = 693.5294
$ws.Range("K107").Value = 6253.5
$ws.Range("L107").Value = 2080.5882
$ws.Range("M107").Value = -4333.5
$ws.Range("N107").Value = -5920.5882
$ws.Range("H113").Value = 554.2857
$ws.Range("I113").Value = 360
$ws.Range("J113").Value = 586.6667
$ws.Range("K113").Value = 1080
$ws.Range("L113").Value = 1760.0001
$ws.Range("M113").Value = 1090
$ws.Range("N113").Value = -6100.0001
$ws.Range("H129").Value = 3600.5334
$ws.Range("J129").Value = 4132
$ws.Range("L129").Value = 12396
$ws.Range("N129").Value = -22396
$ws.Range("H131").Value = 1723.84
$ws.Range("I131").Value = 1130.75
$ws.Range("J131").Value = 2002.9412
$ws.Range("K131").Value = 3392.25
$ws.Range("L131").Value = 6008.8236
$ws.Range("M131").Value = 1647.75
$ws.Range("N131").Value = -16088.8236
$ws.Range("H132").Value = 2187.3
$ws.Range("I132").Value = 2426.5715
$ws.Range("J132").Value = 1629
$ws.Range("K132").Value = 21839.1435
$ws.Range("L132").Value = 14661
$ws.Range("M132").Value = -19309.1435
$ws.Range("N132").Value = -19721

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7373.758
$ws.Range("I70").Value = 7176.5454
$ws.Range("K70").Value = 7176.5454
$ws.Range("M70").Value = -6906.5454
$ws.Range("H73").Value = 7373.758
$ws.Range("I73").Value = 7176.5454
$ws.Range("K73").Value = 7176.5454
$ws.Range("M73").Value = -6240.5454
$ws.Range("H97").Value = 2683.3333
$ws.Range("I97").Value = 1961
$ws.Range("K97").Value = 1961
$ws.Range("M97").Value = -1465
$ws.Range("H132").Value = 3042.04
$ws.Range("I132").Value = 3111.8635
$ws.Range("J132").Value = 2530
$ws.Range("K132").Value = 9335.5905
$ws.Range("L132").Value = 7590
$ws.Range("M132").Value = -6805.5905
$ws.Range("N132").Value = -12650
$ws.Range("H136").Value = 18806.5
$ws.Range("J136").Value = 18806.5
$ws.Range("L136").Value = 56419.5
$ws.Range("N136").Value = -61519.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3100.647
$ws.Range("J46").Value = 4711.2
$ws.Range("L46").Value = 4711.2
$ws.Range("N46").Value = -5087.2
$ws.Range("H132").Value = 2924.3
$ws.Range("I132").Value = 2689.72
$ws.Range("K132").Value = 8069.16
$ws.Range("M132").Value = -5539.16
$ws.Range("H141").Value = 106657.336
$ws.Range("J141").Value = 106657.336
$ws.Range("L141").Value = 106657.336
$ws.Range("N141").Value = -117017.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1506.2858
$ws.Range("I100").Value = 250
$ws.Range("J100").Value = 1715.6666
$ws.Range("K100").Value = 500
$ws.Range("L100").Value = 3431.3332
$ws.Range("M100").Value = 41
$ws.Range("N100").Value = -4513.3332
$ws.Range("H132").Value = 1267.2307
$ws.Range("I132").Value = 1018.5143
$ws.Range("J132").Value = 3443.5
$ws.Range("K132").Value = 3055.5429
$ws.Range("L132").Value = 10330.5
$ws.Range("M132").Value = -525.5429000000004
$ws.Range("N132").Value = -15390.5
